# Apply "added Sea - MGO/HFO; updated accordingly" edit to init_fuel_mix sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("init_fuel_mix")
$ws.Activate()

# Row 9 (Sea / Diesel) -> Sea / MGO
$ws.Range("C9").Value = "MGO"

# Row 10 (Sea / Hydrogen) -> Sea / HFO
$ws.Range("C10").Value = "HFO"

# Insert a new row before the old "Sea / Ammonia" row (row 11), shifting
# Ammonia (old row 11) and Methanol (old row 12) down by one row.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with "Sea / Hydrogen / 0"
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Sea"
$ws.Range("C11").Value = "Hydrogen"
$ws.Range("D11").Value = 0

# Renumber the Index column for the rows that shifted down
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
